# Insert a new price-observation row at row 250 (pushing the existing
# rows 250-360 down to 251-361) and populate the new row with its data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(250).EntireRow.Insert() | Out-Null

$ws.Cells.Item(250, 1).Value2 = 10
$ws.Cells.Item(250, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(250, 3).Value2 = "La Araucanía"
$ws.Cells.Item(250, 4).Value2 = 44704
$ws.Cells.Item(250, 5).Value2 = 9
$ws.Cells.Item(250, 6).Value2 = "Fruta"
$ws.Cells.Item(250, 7).Value2 = 100108
$ws.Cells.Item(250, 8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(250, 9).Value2 = 100108002
$ws.Cells.Item(250, 10).Value2 = "Mango"
$ws.Cells.Item(250, 11).Value2 = "Sin especificar"
$ws.Cells.Item(250, 12).Value2 = "Primera"
$ws.Cells.Item(250, 13).Value2 = 250
$ws.Cells.Item(250, 14).Value2 = 11000
$ws.Cells.Item(250, 15).Value2 = 11000
$ws.Cells.Item(250, 16).Value2 = 11000
$ws.Cells.Item(250, 17).Value2 = "`$/bandeja 4 kilos"
$ws.Cells.Item(250, 18).Value2 = "Perú"
$ws.Cells.Item(250, 19).Value2 = 2750
$ws.Cells.Item(250, 20).Value2 = 4
